$d = $word.ActiveDocument

# Update the header date line (unique text, safe to use Find/Replace)
$d.Content.Find.Execute("2024-03-12 Tuesday", $true, $false, $false, $false, $false, $true, 1, $false, "2024-03-13 Wednesday", 2)

# Update each division-problem cell in the practice table by exact
# (row, column) address -- several old values repeat verbatim in the
# table (e.g. "18÷3=6, 0" and "59÷3=19, 2" each appear twice) but map
# to different replacements, so a global Find/Replace would be unsafe;
# addressing cells directly keeps each substitution independent.
$t = $d.Tables(1)
$t.Cell(1,1).Range.Text = "78÷9=8, 6"
$t.Cell(1,2).Range.Text = "97÷2=48, 1"
$t.Cell(1,3).Range.Text = "26÷3=8, 2"
$t.Cell(1,4).Range.Text = "85÷8=10, 5"
$t.Cell(1,5).Range.Text = "11÷7=1, 4"
$t.Cell(5,1).Range.Text = "76÷5=15, 1"
$t.Cell(5,2).Range.Text = "10÷8=1, 2"
$t.Cell(5,3).Range.Text = "89÷6=14, 5"
$t.Cell(5,4).Range.Text = "26÷6=4, 2"
$t.Cell(5,5).Range.Text = "75÷5=15, 0"
$t.Cell(9,1).Range.Text = "15÷6=2, 3"
$t.Cell(9,2).Range.Text = "27÷2=13, 1"
$t.Cell(9,3).Range.Text = "43÷3=14, 1"
$t.Cell(9,4).Range.Text = "25÷3=8, 1"
$t.Cell(9,5).Range.Text = "55÷4=13, 3"
$t.Cell(13,1).Range.Text = "82÷9=9, 1"
$t.Cell(13,2).Range.Text = "99÷9=11, 0"
$t.Cell(13,3).Range.Text = "14÷4=3, 2"
$t.Cell(13,4).Range.Text = "82÷4=20, 2"
$t.Cell(13,5).Range.Text = "74÷7=10, 4"
$t.Cell(17,1).Range.Text = "45÷9=5, 0"
$t.Cell(17,2).Range.Text = "93÷4=23, 1"
$t.Cell(17,3).Range.Text = "73÷8=9, 1"
$t.Cell(17,4).Range.Text = "51÷2=25, 1"
$t.Cell(17,5).Range.Text = "80÷7=11, 3"
